# msz - 3./4. smoke test + inheritance page and 2. dialog
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename smoke-test identifiers: "VehicleInsuranceAutomobile" -> "AutomobileInsurance"
$ws.Range("A2").Value = "102_AutomobileInsurance_001_SmokeTest"
$ws.Range("B2").Value = "102_AutomobileInsurance_001_SmokeTest_FillPageVehicleData"
$ws.Range("C2").Value = "102_AutomobileInsurance_001_SmokeTest_FillPageInsurantData"
$ws.Range("D2").Value = "102_AutomobileInsurance_001_SmokeTest_FillPageProductData"
$ws.Range("E2").Value = "Choose Platinum"
$ws.Range("F2").Value = "102_AutomobileInsurance_001_SmokeTest_FillPageSendQuote"
$ws.Range("G2").Value = "Send Quote - Button Main Page"

# Move the window on screen and resize it
$excel.ActiveWindow.Left = 10320
$excel.ActiveWindow.Top = 2316
$excel.ActiveWindow.Width = 29856
$excel.ActiveWindow.Height = 12480

# Scroll the view so column C is the left-most visible column, then select D20
$excel.ActiveWindow.ScrollColumn = 3
$ws.Range("D20").Select()
